$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Cálculos")
$win = $excel.ActiveWindow
Write-Host "Active sheet before: $($win.ActiveSheet.Name)"
$ws1.Activate()
Write-Host "Active sheet after activate: $($win.ActiveSheet.Name)"
$win.ScrollRow = 37
Write-Host "ScrollRow: $($win.ScrollRow)"
$ws1.Range("O29").Select()
